# Add data for 2021-10-17:
#  - roll the "current through" date from Oct 08 -> Oct 09 (sheet name + header label)
#  - bump several neighborhood/month counts for the new day's carjacking records

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet and update the running "through" label (shared string) ---
$ws.Name = "Through 2021-10-09"
$ws.Range("B1").Value = "October 2021 (through October 09)"

# --- Row 3: Austin ---
$ws.Range("B3").Value = 4
$ws.Range("L3").Value = 4

# --- Row 4: North Lawndale ---
$ws.Range("L4").Value = 4
$ws.Range("AP4").Value = 1
$ws.Range("AZ4").Value = 1

# --- Row 8: Humboldt Park ---
$ws.Range("B8").Value = 3

# --- Row 9: Grand Crossing ---
$ws.Range("AP9").Value = 2

# --- Row 10: Roseland ---
$ws.Range("AP10").Value = 1

# --- Row 27: West Pullman ---
$ws.Range("BJ27").Value = 1

# --- Row 30: South Chicago ---
$ws.Range("B30").Value = 2

# --- Row 31: Lincoln Park ---
$ws.Range("AF31").Value = 1

# --- Row 37: Englewood ---
$ws.Range("AF37").Value = 1

# --- Row 46: West Elsdon ---
$ws.Range("B46").Value = 1

# --- Row 52: Streeterville ---
$ws.Range("AF52").Value = 1

# --- Row 79: Loop ---
$ws.Range("AP79").Value = 2

# --- Row 87: North Center ---
$ws.Range("AF87").Value = 2

# --- Row 97: South Deering ---
$ws.Range("L97").Value = 1
